$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume/rank data (scrape refresh, GitHub Actions run)

# --- Price column (D): force text format so values like "1.00" or "0.120"
# keep their exact textual representation instead of being coerced to numbers ---
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D11","D12","D13","D14","D15","D16","D17","D18","D21","D23","D24","D25","D27","D28","D29","D32","D33","D35","D37","D39","D40","D41","D43","D44","D45","D50")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$updates = @{
    'D2' = '69.386.26'
    'E2' = '  +2.51%  '
    'D3' = '3.389.53'
    'E3' = '  +1.79%  '
    'D4' = '1.00'
    'E4' = '  +0.08%  '
    'D5' = '589.78'
    'E5' = '  +1.61%  '
    'D6' = '180.80'
    'E6' = '  +2.87%  '
    'D7' = '1.00'
    'E7' = '  +0.03%  '
    'D8' = '0.595'
    'E8' = '  +1.01%  '
    'E9' = '  +8.03%  '
    'E10' = '  +1.45%  '
    'D11' = '48.80'
    'E11' = '  +4.84%  '
    'D12' = '0.0000286'
    'E12' = '  +5.25%  '
    'D13' = '687.34'
    'E13' = '  -2.46%  '
    'D14' = '8.64'
    'E14' = '  +2.29%  '
    'D15' = '3.939.27'
    'E15' = '  +1.61%  '
    'D16' = '69.440.56'
    'E16' = '  +2.63%  '
    'B17' = 'WrappedEther'
    'C17' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D17' = '3.386.61'
    'E17' = '  +1.46%  '
    'B18' = 'TRON'
    'C18' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'D18' = '0.120'
    'E18' = '  +1.75%  '
    'E19' = '  +2.36%  '
    'E20' = '  +3.90%  '
    'D21' = '0.902'
    'E21' = '  +0.91%  '
    'E22' = '  -0.34%  '
    'D23' = '17.06'
    'E23' = '  +0.79%  '
    'D24' = '104.42'
    'E24' = '  +6.29%  '
    'D25' = '3.95'
    'E25' = '  +1.79%  '
    'E26' = '  +1.66%  '
    'D27' = '9.63'
    'E27' = '  +1.19%  '
    'D28' = '34.55'
    'E28' = '  +4.10%  '
    'D29' = '8.70'
    'E29' = '  +1.99%  '
    'E30' = '  -1.31%  '
    'E31' = '  +2.13%  '
    'D32' = '558.34'
    'E32' = '  -2.18%  '
    'D33' = '3.64'
    'E33' = '  +9.38%  '
    'E34' = '  +1.09%  '
    'D35' = '58.19'
    'E35' = '  +1.56%  '
    'E36' = '  +0.22%  '
    'D37' = '3.714.32'
    'E37' = '  +0.28%  '
    'E38' = '  +8.38%  '
    'D39' = '35.11'
    'E39' = '  +3.22%  '
    'D40' = '3.25'
    'E40' = '  +1.66%  '
    'D41' = '0.0₃0707'
    'E41' = '  +4.99%  '
    'E42' = '  +1.50%  '
    'D43' = '0.341'
    'E43' = '  +1.36%  '
    'D44' = '0.0419'
    'E44' = '  +3.18%  '
    'D45' = '3.28'
    'E45' = '  -0.49%  '
    'E46' = '  -0.72%  '
    'E47' = '  +1.34%  '
    'E48' = '  +5.67%  '
    'E49' = '  -0.05%  '
    'D50' = '132.76'
    'E50' = '  +3.30%  '
    'E51' = '  -1.59%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# Restore default (Normal) style on the price cells so no stray number-format
# style index is left attached to them
foreach ($ref in $priceCells) {
    $ws.Range($ref).Style = "Normal"
}

